$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 386.44
$ws.Range("I17").Value = 118.4
$ws.Range("K17").Value = 355.2
$ws.Range("M17").Value = -187.2
$ws.Range("H19").Value = 810
$ws.Range("I19").Value = 640.0714
$ws.Range("J19").Value = 942.1667
$ws.Range("K19").Value = 640.0714
$ws.Range("L19").Value = 942.1667
$ws.Range("M19").Value = -465.0714
$ws.Range("N19").Value = -1292.1667
$ws.Range("H62").Value = 4080.7083
$ws.Range("I62").Value = 1635.4
$ws.Range("J62").Value = 4724.2104
$ws.Range("K62").Value = 1635.4
$ws.Range("L62").Value = 4724.2104
$ws.Range("M62").Value = -1011.4
$ws.Range("N62").Value = -5972.2104
$ws.Range("H65").Value = 4080.7083
$ws.Range("I65").Value = 1635.4
$ws.Range("J65").Value = 4724.2104
$ws.Range("K65").Value = 8177
$ws.Range("L65").Value = 23621.052
$ws.Range("M65").Value = -5057
$ws.Range("N65").Value = -29861.052
$ws.Range("H115").Value = 436.125
$ws.Range("I115").Value = 436.125
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1308.375
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 258.625
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 1136.4706
$ws.Range("I118").Value = 665
$ws.Range("J118").Value = 1555.5555
$ws.Range("K118").Value = 1995
$ws.Range("L118").Value = 4666.666499999999
$ws.Range("M118").Value = -338
$ws.Range("N118").Value = -7980.666499999999
$ws.Range("H127").Value = 925.3
$ws.Range("I127").Value = 512
$ws.Range("J127").Value = 1263.4546
$ws.Range("K127").Value = 1536
$ws.Range("L127").Value = 3790.3638
$ws.Range("M127").Value = 3424
$ws.Range("N127").Value = -13710.3638
$ws.Range("H129").Value = 1050.6086
$ws.Range("I129").Value = 606.7143
$ws.Range("J129").Value = 1244.8125
$ws.Range("K129").Value = 1820.1429
$ws.Range("L129").Value = 3734.4375
$ws.Range("M129").Value = 3179.8571
$ws.Range("N129").Value = -13734.4375
$ws.Range("H135").Value = 626229.75
$ws.Range("I135").Value = 280.2
$ws.Range("J135").Value = 2311478.5
$ws.Range("K135").Value = 2521.8
$ws.Range("L135").Value = 20803306.5
$ws.Range("M135").Value = 13.20000000000027
$ws.Range("N135").Value = -20808376.5
$ws.Range("H137").Value = 48612600
$ws.Range("I137").Value = 8334678
$ws.Range("J137").Value = 250002200
$ws.Range("K137").Value = 25004034
$ws.Range("L137").Value = 750006600
$ws.Range("M137").Value = -25001484
$ws.Range("N137").Value = -750011700
$ws.Range("H141").Value = 1025.5
$ws.Range("I141").Value = 536.4
$ws.Range("J141").Value = 2073.5715
$ws.Range("K141").Value = 1609.2
$ws.Range("L141").Value = 6220.7145
$ws.Range("M141").Value = 3570.8
$ws.Range("N141").Value = -16580.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2264003.5
$ws.Range("I2").Value = 1682.2
$ws.Range("J2").Value = 3677954.2
$ws.Range("K2").Value = 1682.2
$ws.Range("L2").Value = 3677954.2
$ws.Range("M2").Value = -1569.2
$ws.Range("N2").Value = -3678180.2
$ws.Range("H32").Value = 5514.322
$ws.Range("I32").Value = 3361.2253
$ws.Range("J32").Value = 15068.6875
$ws.Range("K32").Value = 3361.2253
$ws.Range("L32").Value = 15068.6875
$ws.Range("M32").Value = -3074.2253
$ws.Range("N32").Value = -15642.6875
$ws.Range("H61").Value = 13335076
$ws.Range("I61").Value = 18520222
$ws.Range("K61").Value = 18520222
$ws.Range("M61").Value = -18520010
$ws.Range("H116").Value = 2264003.5
$ws.Range("I116").Value = 1682.2
$ws.Range("J116").Value = 3677954.2
$ws.Range("K116").Value = 1682.2
$ws.Range("L116").Value = 3677954.2
$ws.Range("M116").Value = 611.8
$ws.Range("N116").Value = -3682542.2
$ws.Range("H136").Value = 13335076
$ws.Range("I136").Value = 18520222
$ws.Range("K136").Value = 55560666
$ws.Range("M136").Value = -55558116

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2264003.5
$ws.Range("I3").Value = 1682.2
$ws.Range("J3").Value = 3677954.2
$ws.Range("K3").Value = 1682.2
$ws.Range("L3").Value = 3677954.2
$ws.Range("M3").Value = -1568.2
$ws.Range("N3").Value = -3678182.2
$ws.Range("H107").Value = 1290.25
$ws.Range("I107").Value = 1231.6923
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 1231.6923
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = 688.3077000000001
$ws.Range("N107").Value = -5239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2290.0715
$ws.Range("I16").Value = 3388.875
$ws.Range("J16").Value = 825
$ws.Range("K16").Value = 3388.875
$ws.Range("L16").Value = 825
$ws.Range("M16").Value = -3101.875
$ws.Range("N16").Value = -1399
$ws.Range("H31").Value = 1653.4359
$ws.Range("I31").Value = 1357.5
$ws.Range("J31").Value = 1819.16
$ws.Range("K31").Value = 1357.5
$ws.Range("L31").Value = 1819.16
$ws.Range("M31").Value = -1062.5
$ws.Range("N31").Value = -2409.16
$ws.Range("H34").Value = 1653.4359
$ws.Range("I34").Value = 1357.5
$ws.Range("J34").Value = 1819.16
$ws.Range("K34").Value = 1357.5
$ws.Range("L34").Value = 1819.16
$ws.Range("M34").Value = -1155.5
$ws.Range("N34").Value = -2223.16
$ws.Range("H94").Value = 787.25
$ws.Range("I94").Value = 617.8333
$ws.Range("J94").Value = 859.8570999999999
$ws.Range("K94").Value = 617.8333
$ws.Range("L94").Value = 859.8570999999999
$ws.Range("M94").Value = -166.8333
$ws.Range("N94").Value = -1761.8571
$ws.Range("H113").Value = 2290.0715
$ws.Range("I113").Value = 3388.875
$ws.Range("J113").Value = 825
$ws.Range("K113").Value = 3388.875
$ws.Range("L113").Value = 825
$ws.Range("M113").Value = -1218.875
$ws.Range("N113").Value = -5165
$ws.Range("H132").Value = 3169
$ws.Range("I132").Value = 2947.2778
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 8841.8334
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -6311.8334
$ws.Range("N132").Value = -18558.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 612.5
$ws.Range("J113").Value = 600
$ws.Range("L113").Value = 1800
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 27780506
$ws.Range("I122").Value = 66668064
$ws.Range("J122").Value = 3678.4285
$ws.Range("K122").Value = 600012576
$ws.Range("L122").Value = 33105.8565
$ws.Range("M122").Value = -600010126
$ws.Range("N122").Value = -38005.8565
$ws.Range("H132").Value = 333335330
$ws.Range("I132").Value = 1000000000
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9000000000
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -8999997470
$ws.Range("N132").Value = -32060
$ws.Range("H133").Value = 6167.1333
$ws.Range("I133").Value = 3463.375
$ws.Range("J133").Value = 9257.143
$ws.Range("K133").Value = 10390.125
$ws.Range("L133").Value = 27771.429
$ws.Range("M133").Value = -5330.125
$ws.Range("N133").Value = -37891.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 524.34784
$ws.Range("I107").Value = 347
$ws.Range("J107").Value = 638.3570999999999
$ws.Range("K107").Value = 347
$ws.Range("L107").Value = 638.3570999999999
$ws.Range("M107").Value = 1573
$ws.Range("N107").Value = -4478.3571
$ws.Range("H113").Value = 5319.0713
$ws.Range("I113").Value = 7959.1875
$ws.Range("J113").Value = 1798.9166
$ws.Range("K113").Value = 7959.1875
$ws.Range("L113").Value = 1798.9166
$ws.Range("M113").Value = -5789.1875
$ws.Range("N113").Value = -6138.9166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3450
$ws.Range("I7").Value = 3466.6667
$ws.Range("K7").Value = 3466.6667
$ws.Range("M7").Value = -3354.6667
$ws.Range("H93").Value = 2844.7778
$ws.Range("I93").Value = 2929
$ws.Range("J93").Value = 2550
$ws.Range("K93").Value = 2929
$ws.Range("L93").Value = 2550
$ws.Range("M93").Value = -1681
$ws.Range("N93").Value = -5046
$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -12250
$ws.Range("H126").Value = 3450
$ws.Range("I126").Value = 3466.6667
$ws.Range("K126").Value = 10400.0001
$ws.Range("M126").Value = -7930.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5691.6763
$ws.Range("I132").Value = 5807.0645
$ws.Range("K132").Value = 17421.1935
$ws.Range("M132").Value = -14891.1935
$ws.Range("H136").Value = 7187.375
$ws.Range("I136").Value = 8456.23
$ws.Range("J136").Value = 1689
$ws.Range("K136").Value = 25368.69
$ws.Range("L136").Value = 5067
$ws.Range("M136").Value = -22818.69
$ws.Range("N136").Value = -10167
